$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.778.06"
$ws.Cells.Item(2, 5).Value = "  +1.17%  "

$ws.Cells.Item(3, 4).Value = "1.854.93"
$ws.Cells.Item(3, 5).Value = "  +0.89%  "

$ws.Cells.Item(5, 4).Value = "'244.08"
$ws.Cells.Item(5, 5).Value = "  +0.31%  "

$ws.Cells.Item(6, 4).Value = "'0.6535"
$ws.Cells.Item(6, 5).Value = "  +4.12%  "

$ws.Cells.Item(7, 5).Value = "  +0.17%  "

$ws.Cells.Item(8, 4).Value = "'47.97"
$ws.Cells.Item(8, 5).Value = "  +4.06%  "

$ws.Cells.Item(9, 4).Value = "'0.07540"
$ws.Cells.Item(9, 5).Value = "  +1.66%  "

$ws.Cells.Item(10, 4).Value = "'0.2967"
$ws.Cells.Item(10, 5).Value = "  +0.63%  "

$ws.Cells.Item(11, 4).Value = "'24.58"
$ws.Cells.Item(11, 5).Value = "  +4.93%  "

$ws.Cells.Item(12, 4).Value = "'0.07638"
$ws.Cells.Item(12, 5).Value = "  -0.11%  "

$ws.Cells.Item(13, 4).Value = "1.858.60"
$ws.Cells.Item(13, 5).Value = "  +1.07%  "

$ws.Cells.Item(14, 4).Value = "'5.058"
$ws.Cells.Item(14, 5).Value = "  +0.87%  "

$ws.Cells.Item(15, 4).Value = "'0.6877"
$ws.Cells.Item(15, 5).Value = "  +1.74%  "

$ws.Cells.Item(16, 4).Value = "'83.70"
$ws.Cells.Item(16, 5).Value = "  +0.34%  "

$ws.Cells.Item(17, 4).Value = "'0.000009719"
$ws.Cells.Item(17, 5).Value = "  +3.93%  "

$ws.Cells.Item(18, 4).Value = "'6.110"
$ws.Cells.Item(18, 5).Value = "  +3.36%  "

$ws.Cells.Item(19, 4).Value = "29.807.84"
$ws.Cells.Item(19, 5).Value = "  +1.34%  "

$ws.Cells.Item(20, 4).Value = "2.113.46"
$ws.Cells.Item(20, 5).Value = "  +1.48%  "

$ws.Cells.Item(21, 4).Value = "'236.43"
$ws.Cells.Item(21, 5).Value = "  -0.48%  "

$ws.Cells.Item(22, 4).Value = "'12.65"
$ws.Cells.Item(22, 5).Value = "  +0.80%  "

$ws.Cells.Item(23, 2).Value = "Dai"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(23, 4).Value = "'1.001"
$ws.Cells.Item(23, 5).Value = "  +0.13%  "

$ws.Cells.Item(24, 2).Value = "Chainlink"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(24, 4).Value = "'7.771"
$ws.Cells.Item(24, 5).Value = "  +5.97%  "

$ws.Cells.Item(25, 5).Value = "  +0.25%  "

$ws.Cells.Item(26, 4).Value = "'158.37"
$ws.Cells.Item(26, 5).Value = "  -0.39%  "

$ws.Cells.Item(27, 4).Value = "'0.1433"
$ws.Cells.Item(27, 5).Value = "  +1.31%  "

$ws.Cells.Item(28, 4).Value = "'8.556"
$ws.Cells.Item(28, 5).Value = "  +0.61%  "

$ws.Cells.Item(29, 4).Value = "'17.88"
$ws.Cells.Item(29, 5).Value = "  +0.70%  "

$ws.Cells.Item(30, 5).Value = "  -0.11%  "

$ws.Cells.Item(31, 4).Value = "'0.06036"
$ws.Cells.Item(31, 5).Value = "  +0.44%  "

$ws.Cells.Item(32, 4).Value = "'1.278"
$ws.Cells.Item(32, 5).Value = "  +3.57%  "

$ws.Cells.Item(33, 4).Value = "'4.154"
$ws.Cells.Item(33, 5).Value = "  +1.09%  "

$ws.Cells.Item(34, 4).Value = "'4.091"
$ws.Cells.Item(34, 5).Value = "  -0.07%  "

$ws.Cells.Item(35, 5).Value = "  +0.26%  "

$ws.Cells.Item(36, 4).Value = "'1.179"
$ws.Cells.Item(36, 5).Value = "  +3.38%  "

$ws.Cells.Item(37, 4).Value = "'0.7255"
$ws.Cells.Item(37, 5).Value = "  -0.03%  "

$ws.Cells.Item(38, 4).Value = "'2.611"
$ws.Cells.Item(38, 5).Value = "  +0.03%  "

$ws.Cells.Item(39, 4).Value = "'2.808"
$ws.Cells.Item(39, 5).Value = "  -2.55%  "

$ws.Cells.Item(40, 4).Value = "'0.01792"
$ws.Cells.Item(40, 5).Value = "  +1.97%  "

$ws.Cells.Item(41, 4).Value = "1.204.84"
$ws.Cells.Item(41, 5).Value = "  -1.05%  "

$ws.Cells.Item(42, 4).Value = "'6.301"
$ws.Cells.Item(42, 5).Value = "  +0.38%  "

$ws.Cells.Item(43, 4).Value = "'0.9137"
$ws.Cells.Item(43, 5).Value = "  +0.24%  "

$ws.Cells.Item(44, 4).Value = "'1.001"
$ws.Cells.Item(44, 5).Value = "  +0.03%  "

$ws.Cells.Item(45, 4).Value = "2.022.17"
$ws.Cells.Item(45, 5).Value = "  +1.25%  "

$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "'67.13"
$ws.Cells.Item(46, 5).Value = "  +2.31%  "

$ws.Cells.Item(47, 2).Value = "Quant"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(47, 4).Value = "'101.15"
$ws.Cells.Item(47, 5).Value = "  -0.69%  "

$ws.Cells.Item(48, 5).Value = "  +1.30%  "

$ws.Cells.Item(49, 4).Value = "'7.291"
$ws.Cells.Item(49, 5).Value = "  +9.48%  "

$ws.Cells.Item(50, 4).Value = "'0.4059"
$ws.Cells.Item(50, 5).Value = "  -0.02%  "

$ws.Cells.Item(51, 4).Value = "'9.130"
$ws.Cells.Item(51, 5).Value = "  -1.01%  "
